$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Hunk 1 — paragraph 2 ("Uma casa ...")
# Split " está associada a um ou mais utilizadores. Cada casa é
# caracterizada por um identificador único" into three sentences by
# inserting ", podendo um utilizador ter várias casas" before the
# final sentence.
# -----------------------------------------------------------------
$needle1 = " está associada a um ou mais utilizadores. Cada casa é caracterizada por um identificador único"
$replace1 = " está associada a um ou mais utilizadores, podendo um utilizador ter várias casas. Cada casa é caracterizada por um identificador único"
$p2 = $d.Paragraphs.Item(2).Range
$found1 = $p2.Find.Execute($needle1, $true, $false, $false, $false, $false, $true, 1, $false, $replace1, 2)
Write-Host "Hunk1 found: $found1"

# -----------------------------------------------------------------
# Hunk 2 — paragraph 5 ("Um item presente numa casa ...")
# Rewrite the item-characterisation sentence, dropping the
# sku/(stock keeping unit) mention, the yellow highlight and the
# italics, and moving "uma descrição" earlier in the sentence.
# -----------------------------------------------------------------
$needle2 = "é designado por um identificador, um sku (stock keeping unit), uma variedade, uma marca, um segmento, e a sua unidade de medida (por exemplo, litro, mililitro, etc.)"
$replace2 = "é identificado por um identificador único ou por uma marca, uma variedade e um segmento, é também caracterizado por uma descrição"
$text = $d.Content.Text
$idx2 = $text.IndexOf($needle2)
if ($idx2 -lt 0) { Write-Host "Hunk2a NOT FOUND" }
$r2 = $d.Range($idx2, $idx2 + $needle2.Length)
$r2.Text = $replace2

# Remove the now-redundant ", uma descrição" leftover (old copy of
# "uma descrição" that used to follow the comma).
$text = $d.Content.Text
$needle2b = ", uma descrição o local de conservação"
$idx2b = $text.IndexOf($needle2b)
if ($idx2b -lt 0) { Write-Host "Hunk2b NOT FOUND" }
$delStart = $idx2b + 1
$delLen = " uma descrição".Length
$r2b = $d.Range($delStart, $delStart + $delLen)
$r2b.Text = ""

# -----------------------------------------------------------------
# Hunk 3 — still paragraph 5, the movements sentence.
# Remove the cyan highlight from "de um", drop the stray space run
# that followed it, and shrink "(do seu) " (red) down to a lone
# space while keeping the red colour.
# -----------------------------------------------------------------
$text = $d.Content.Text
$needle3 = "saiu de um (do seu) local"
$idx3 = $text.IndexOf($needle3)
if ($idx3 -lt 0) { Write-Host "Hunk3a NOT FOUND" }
$startDeUm = $idx3 + "saiu ".Length
$rClear = $d.Range($startDeUm, $startDeUm + "de um ".Length)
$rClear.Text = ""
$rIns = $d.Range($startDeUm, $startDeUm)
$rIns.InsertAfter("de um")

$text = $d.Content.Text
$needle3b = "(do seu) "
$idx3b = $text.IndexOf($needle3b)
if ($idx3b -lt 0) { Write-Host "Hunk3b NOT FOUND" }
$r3b = $d.Range($idx3b, $idx3b + $needle3b.Length)
$r3b.Text = " "

# -----------------------------------------------------------------
# Hunk 3 (cont.) — drop ", o local de armazenamento envolvido" and
# move the separating space so it sits right before the _GoBack
# bookmark instead of right after it.
# -----------------------------------------------------------------
$text = $d.Content.Text
$needle3c = ", o local de armazenamento envolvido"
$idx3c = $text.IndexOf($needle3c)
if ($idx3c -lt 0) { Write-Host "Hunk3c NOT FOUND" }
$r3c = $d.Range($idx3c, $idx3c + $needle3c.Length)
$r3c.Text = ""

$text = $d.Content.Text
$needle3d = " e a quantidade de produtos que ocorrem num movimento."
$idx3d = $text.IndexOf($needle3d)
if ($idx3d -lt 0) { Write-Host "Hunk3d NOT FOUND" }
$r3d = $d.Range($idx3d, $idx3d + 1)
$r3d.Text = ""

$text = $d.Content.Text
$needle3e = "a data em que ocorreu o movimento"
$idx3e = $text.IndexOf($needle3e)
if ($idx3e -lt 0) { Write-Host "Hunk3e NOT FOUND" }
$insertAt = $idx3e + $needle3e.Length
$r3e = $d.Range($insertAt, $insertAt)
$r3e.InsertAfter(" ")

Write-Host "Done"
